$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.522.60'
$ws.Range("E2").Value = '  -0.65%  '
$ws.Range("D3").Value = '2.081.01'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.632'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.39%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.84'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0781'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.55%  '
$ws.Range("E11").Value = '  +2.43%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.07'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("D13").Value = '2.387.63'
$ws.Range("E13").Value = '  +0.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.07'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.773'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.35'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.20%  '
$ws.Range("D17").Value = '2.077.11'
$ws.Range("E17").Value = '  -0.30%  '
$ws.Range("D18").Value = '37.465.65'
$ws.Range("E18").Value = '  -0.57%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.05'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.24%  '
$ws.Range("B20").Value = 'Litecoin'
$ws.Range("C20").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.75'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '0.0₃0833'
$ws.Range("E21").Value = '  -0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '228.70'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.10%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.34'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.36%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("E28").Value = '  -4.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.56'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("E30").Value = '  -0.92%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.123'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.74%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0636'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.67'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.48'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("E36").Value = '  -0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.32'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.29'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.28%  '
$ws.Range("E40").Value = '  +7.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.30'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0956'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.20'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.92'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.82%  '
$ws.Range("D46").Value = '1.457.59'
$ws.Range("E46").Value = '  +0.40%  '
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("B48").Value = 'FraxShare'
$ws.Range("C48").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.94'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.44%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.95'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.98%  '
$ws.Range("D51").Value = '2.271.09'
$ws.Range("E51").Value = '  +0.04%  '
